$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row-1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 data values (B2:E2)
$ws.Range("B2").Value = 19.784121002394567
$ws.Range("C2").Value = 42.012862100795893
$ws.Range("D2").Value = 18.762756017546945
$ws.Range("E2").Value = 37.510580235727396

# Update row 3 data values (B3:E3)
$ws.Range("B3").Value = 19.421284614683866
$ws.Range("C3").Value = 30.004742429034
$ws.Range("D3").Value = 24.705470356675846
$ws.Range("E3").Value = 23.702137112784595

# Update the sheet selection to match the new, smaller edited range
$ws.Range("B1:E3").Select()
